# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.758.58"
$ws.Range("E2").Value = "  -1.79%  "

# Row 3
$ws.Range("D3").Value = "2.301.37"
$ws.Range("E3").Value = "  -2.82%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.36"
$ws.Range("E5").Value = "  -1.15%  "

# Row 6
$ws.Range("E6").Value = "  -3.25%  "

# Row 7
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("E8").Value = "  -3.25%  "

# Row 9
$ws.Range("D9").Value = "2.300.27"
$ws.Range("E9").Value = "  -2.76%  "

# Row 10
$ws.Range("E10").Value = "  -1.01%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.50"
$ws.Range("E11").Value = "  -0.14%  "

# Row 12
$ws.Range("E12").Value = "  -0.69%  "

# Row 13
$ws.Range("E13").Value = "  -2.11%  "

# Row 14
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "59.697.51"
$ws.Range("E14").Value = "  -1.74%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.708.86"
$ws.Range("E15").Value = "  -3.03%  "

# Row 16
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.04"
$ws.Range("E16").Value = "  -4.81%  "

# Row 17
$ws.Range("E17").Value = "  -2.11%  "

# Row 18
$ws.Range("D18").Value = "2.305.57"
$ws.Range("E18").Value = "  -3.57%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.39"
$ws.Range("E19").Value = "  -3.21%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.02"
$ws.Range("E20").Value = "  -4.56%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "309.45"
$ws.Range("E21").Value = "  -2.58%  "

# Row 22
$ws.Range("E22").Value = "  -8.36%  "

# Row 23
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.07"
$ws.Range("E24").Value = "  -0.61%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("E25").Value = "  -2.53%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.14%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.69"
$ws.Range("E27").Value = "  -4.94%  "

# Row 28
$ws.Range("E28").Value = "  -0.69%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.34"
$ws.Range("E29").Value = "  +0.22%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").Value = "  +3.90%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.70"
$ws.Range("E31").Value = "  -2.44%  "

# Row 32
$ws.Range("D32").Value = "0.0₃0710"
$ws.Range("E32").Value = "  -5.11%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.77"
$ws.Range("E33").Value = "  -2.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.377"
$ws.Range("E34").Value = "  -1.89%  "

# Row 35
$ws.Range("E35").Value = "  +0.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.33"
$ws.Range("E36").Value = "  -5.54%  "

# Row 37
$ws.Range("E37").Value = "  -2.58%  "

# Row 38
$ws.Range("E38").Value = "  +0.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.99"
$ws.Range("E39").Value = "  -5.76%  "

# Row 40
$ws.Range("E40").Value = "  -3.76%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.76"
$ws.Range("E41").Value = "  -1.47%  "

# Row 42
$ws.Range("E42").Value = "  -4.36%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.37"
$ws.Range("E43").Value = "  -6.21%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.40"
$ws.Range("E44").Value = "  -2.53%  "

# Row 45
$ws.Range("E45").Value = "  -2.08%  "

# Row 46
$ws.Range("E46").Value = "  -0.29%  "

# Row 47
$ws.Range("E47").Value = "  -3.05%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.42"
$ws.Range("E48").Value = "  -6.01%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0222"
$ws.Range("E49").Value = "  +8.56%  "

# Row 50
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0211"
$ws.Range("E50").Value = "  -1.31%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.99"
$ws.Range("E51").Value = "  -0.44%  "
